# The SMA antenna connector (SMA-KWE903 / RF1, BOM item "7") was dropped
# from this PCB revision, so its row is removed from the BOM and every row
# below it shifts up by one. The rows that move up keep their original
# "No." labels (7, 8), and the worksheet/tab name switches from the SMA
# variant to the NO_ANT (no antenna) variant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the E22-900M30S (was row 9) and AMS1117-3.3V (was row 10) rows up
# onto row 8 (the now-removed SMA-KWE903 row) and row 9, via copy/paste so
# their original cell formatting/types carry over faithfully.
$ws.Range("A9:J10").Copy()
$ws.Range("A8:J9").PasteSpecial()

# Drop the now-duplicated old row 10; the trailing blank row shifts up to
# become row 10.
$ws.Rows.Item(10).Delete()

# Renumber the "No." column for the two rows that moved up (text, like the
# rest of that column - force it past Excel's automatic number detection).
$ws.Cells.Item(8, 1).Value2 = "'7"
$ws.Cells.Item(9, 1).Value2 = "'8"

# Reflect the antenna-less revision in the sheet/tab name.
$ws.Name = "BOM_PCB V1.1_NO_ANT_2024-09-07"
